$wb = $excel.ActiveWorkbook

# --- Sheet: weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.83545108599723
$ws.Range("C2").Value = 0.159160874987795
$ws.Range("B3").Value = 0.145673114301433
$ws.Range("C3").Value = 0.0977476357725525

# --- Sheet: lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.25372991308615
$ws.Range("C2").Value = 0.193599458785533
$ws.Range("B3").Value = -1.05166003595957
$ws.Range("C3").Value = 0.0926540495239974

# --- Sheet: llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.14081736529117
$ws.Range("C2").Value = 0.0953969328983439
$ws.Range("B3").Value = 1.80317848000473
$ws.Range("C3").Value = 0.177904693938904

# --- Sheet: gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.53105087692342
$ws.Range("C2").Value = 0.130124689890718
$ws.Range("B3").Value = -0.00238000384600223
$ws.Range("C3").Value = 0.0140000910362655

# --- Sheet: weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0253321841268806
$ws.Range("B2").Value = -0.0120621687812353
$ws.Range("A3").Value = -0.0120621687812353
$ws.Range("B3").Value = 0.00955460029912358

# --- Sheet: lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0374807504420513
$ws.Range("B2").Value = -0.0159888332730207
$ws.Range("A3").Value = -0.0159888332730207
$ws.Range("B3").Value = 0.00858477289319537

# --- Sheet: llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00910057480641112
$ws.Range("B2").Value = 0.00299024922680052
$ws.Range("A3").Value = 0.00299024922680052
$ws.Range("B3").Value = 0.0316500801254951

# --- Sheet: gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0169324349191557
$ws.Range("B2").Value = -0.00125060303915194
$ws.Range("A3").Value = -0.00125060303915194
$ws.Range("B3").Value = 0.000196002549023721
